$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the style (s="1") from the
# existing header cell H1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-50: column I (9) and column J (10) values.
$data = @(
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(6,7),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,9),
    @(7,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,10),
    @(9,9),
    @(6,6),
    @(7,9),
    @(7,7),
    @(5,6),
    @(9,9),
    @(9,9),
    @(5,6),
    @(6,7),
    @(5,6),
    @(7,8),
    @(7,7),
    @(4,6),
    @(8,8),
    @(6,7),
    @(2,4),
    @(6,7),
    @(7,8),
    @(6,7),
    @(5,6),
    @(7,8),
    @(2,4),
    @(6,8),
    @(7,8),
    @(7,7),
    @(9,9),
    @(6,8),
    @(7,8),
    @(4,5),
    @(7,7),
    @(7,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
